$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the standalone "Meta description: ..." paragraph that follows the
#    H1 title at the top of the document.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Blood Moon Wilds Free Today -
#    Exciting Werewolf Slot Game") right before the very last paragraph
#    (the one that used to hold the "Prompt: ..." image-generation text).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

$newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Blood Moon Wilds Free Today - Exciting Werewolf Slot Game</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$lastPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($count)
$newPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 3) Swap the text of the final paragraph (still italic) from the old DALLE
#    image prompt to the new meta-description-style blurb.
# ---------------------------------------------------------------------------
$finalCount = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($finalCount)

$oldPromptText = "Prompt: Create a feature image fitting the game Blood Moon Wilds. DALLE, please create a cartoon-style feature image for Blood Moon Wilds that showcases a happy Maya warrior wearing glasses. The image should incorporate elements of the eerie bayou surrounding New Orleans, such as a full moon shining in the background and werewolves lurking in the shadows. The Maya warrior should be holding a treasure chest filled with gold coins and precious jewels, to represent the potential for big wins in the game. Be creative and use bold, vibrant colors to make the image stand out and capture the attention of online slot players."
$newBlurbText = "Try Blood Moon Wilds slot game for free today and discover exciting werewolf characters, lunar calendar feature, and more. Compatible on all devices."

$finalPara.Range.Find.Execute($oldPromptText, $true, $false, $false, $false, $false, $true, 1, $false, $newBlurbText, 2)

Write-Host "Edit complete. Paragraph count: " $d.Paragraphs.Count
